# Slide 7 (sldId 265 / creationId 628037149): the ROC-curve picture had the
# wrong axis labels baked in; it was swapped out for a corrected export of
# the same chart. Re-create that swap: drop the old embedded picture and
# insert a replacement (same underlying image) repositioned/resized, placed
# after "TextBox 8" in the shape order, matching the author's fix.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Locate the existing picture shape (id=8, "Picture 7") by its Id, rather
# than assuming index order. Fall back to matching by Name in case ids ever
# shift.
$oldPic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 8) { $oldPic = $sh }
}
if ($oldPic -eq $null) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -eq "Picture 7") { $oldPic = $sh }
    }
}

# The shape-id allocator hands out the lowest free id >= 2. Reserve id=2
# with a throwaway shape so that the real replacement picture below lands on
# id=3 (matching the target), then discard the placeholder afterwards.
$reserve = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)

# Duplicate the original picture so the new shape keeps the same embedded
# image (same relationship/bytes) instead of inserting a brand-new media
# file; Duplicate() appends the clone at the end of the shape stack, i.e.
# after "TextBox 8", which is exactly where the corrected picture belongs.
$newPic = $oldPic.Duplicate()

# Remove the old (wrong-labelled) picture and the id=2 placeholder.
$oldPic.Delete()
$reserve.Delete()

# Rename + relabel to match the corrected picture's metadata.
$newPic.Name = "Picture 2"
$newPic.AlternativeText = "A picture containing text, diagram, line, plot`n`nDescription automatically generated"

# Reposition/resize to the corrected picture's frame (values are EMU/12700,
# expressed with enough precision to survive the points round-trip exactly).
$newPic.Left = 9.396535433070866
$newPic.Top = 111.24224409448819
$newPic.Width = 635.0311811023622
$newPic.Height = 341.64574803149605
